$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column K (24-jun) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("J1").Copy()
$ws1.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("K1").Value = "24-jun"

$k2values = @(81.34, 66.17, 63.79, 45.14, 46.62, 71.3, 84.53, 90.09, 73.55, 37.04, 1.34, 0, -0.01, -0.01, -0.03, -0.01, 7.5, 20.06, 71.95, 104.61, 125.4, 125.49, 131.91, 102.48)

for ($i = 0; $i -lt $k2values.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 11).Value = $k2values[$i]
}

# --- Sheet "Gaz": add row 7 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A7").NumberFormat = "@"
$ws2.Range("A7").Value = "2025-06-23"
$ws2.Range("A7").Style = "Normal"
$ws2.Range("B7").Value = 40.9

# --- Sheet "CO2": add row 7 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A7").NumberFormat = "@"
$ws3.Range("A7").Value = "2025-06-23"
$ws3.Range("A7").Style = "Normal"
$ws3.Range("B7").Value = 71.88
